# "returning books works now"
# Book #1 (row 2) has been returned: it is no longer lent out, so clear the
# lending/return info and flip the "Lent" flag off. The reservation-until
# date is cleared too (reservation has lapsed now that the copy is back).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lent -> FALSE
$ws.Range("G2").Value = $false

# Lent to / Lent date / Return date / Reserved until -> blank
# (assign an empty text value via a lone apostrophe so the cell becomes an
# empty *text* cell rather than being removed outright, then drop the
# quote-prefix style that produces so the cell keeps the default style.)
$ws.Range("H2").Value = "'"
$ws.Range("H2").Style = "Normal"

$ws.Range("I2").Value = "'"
$ws.Range("I2").Style = "Normal"

$ws.Range("J2").Value = "'"
$ws.Range("J2").Style = "Normal"

$ws.Range("M2").Value = "'"
$ws.Range("M2").Style = "Normal"
